$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 176.85715
$ws.Range("I11").Value = 176.85715
$ws.Range("K11").Value = 176.85715
$ws.Range("M11").Value = -36.85714999999999

$ws.Range("H96").Value = 698.3333
$ws.Range("I96").Value = 972.5
$ws.Range("J96").Value = 150
$ws.Range("K96").Value = 2917.5
$ws.Range("L96").Value = 450
$ws.Range("M96").Value = -1544.5
$ws.Range("N96").Value = -3196

$ws.Range("H132").Value = 16562.975
$ws.Range("I132").Value = 2960.9644
$ws.Range("K132").Value = 8882.893199999999
$ws.Range("M132").Value = -6352.893199999999

$ws.Range("H138").Value = 3420.7646
$ws.Range("J138").Value = 3939.1516
$ws.Range("L138").Value = 11817.4548
$ws.Range("N138").Value = -22097.4548

$ws.Range("H141").Value = 5214.148
$ws.Range("I141").Value = 4309.227
$ws.Range("K141").Value = 12927.681
$ws.Range("M141").Value = -7747.681

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 18171.477
$ws.Range("I2").Value = 22112.75
$ws.Range("K2").Value = 22112.75
$ws.Range("M2").Value = -21999.75

$ws.Range("H19").Value = 3092
$ws.Range("I19").Value = 1618.25
$ws.Range("K19").Value = 1618.25
$ws.Range("M19").Value = -1389.25

$ws.Range("H26").Value = 1987.5
$ws.Range("I26").Value = 983.3333
$ws.Range("K26").Value = 983.3333
$ws.Range("M26").Value = -653.3333

$ws.Range("H36").Value = 131132.75
$ws.Range("J36").Value = 15000
$ws.Range("L36").Value = 15000
$ws.Range("N36").Value = -15692

$ws.Range("H44").Value = 50000
$ws.Range("J44").Value = 50000
$ws.Range("L44").Value = 50000
$ws.Range("N44").Value = -50976

$ws.Range("H45").Value = 5869.0835
$ws.Range("I45").Value = 5607.5
$ws.Range("J45").Value = 5999.875
$ws.Range("K45").Value = 5607.5
$ws.Range("L45").Value = 5999.875
$ws.Range("M45").Value = -5230.5
$ws.Range("N45").Value = -6753.875

$ws.Range("H55").Value = 33500
$ws.Range("I55").Value = 33500
$ws.Range("K55").Value = 33500
$ws.Range("M55").Value = -33185

$ws.Range("H70").Value = 200577
$ws.Range("J70").Value = 200577
$ws.Range("L70").Value = 200577
$ws.Range("N70").Value = -201117

$ws.Range("H73").Value = 200577
$ws.Range("J73").Value = 200577
$ws.Range("L73").Value = 200577
$ws.Range("N73").Value = -202449

$ws.Range("H102").Value = 2821.25
$ws.Range("I102").Value = 2821.25
$ws.Range("K102").Value = 2821.25
$ws.Range("M102").Value = -1199.25

$ws.Range("H108").Value = 69999
$ws.Range("J108").Value = 69999
$ws.Range("L108").Value = 69999
$ws.Range("N108").Value = -77679

$ws.Range("H116").Value = 18171.477
$ws.Range("I116").Value = 22112.75
$ws.Range("K116").Value = 22112.75
$ws.Range("M116").Value = -19818.75

$ws.Range("H122").Value = 4938.9375
$ws.Range("I122").Value = 4419.4165
$ws.Range("J122").Value = 6497.5
$ws.Range("K122").Value = 13258.2495
$ws.Range("L122").Value = 19492.5
$ws.Range("M122").Value = -10808.2495
$ws.Range("N122").Value = -24392.5

$ws.Range("H132").Value = 2869.0908
$ws.Range("I132").Value = 2531
$ws.Range("K132").Value = 7593
$ws.Range("M132").Value = -5063

$ws.Range("H141").Value = 73000
$ws.Range("J141").Value = 73000
$ws.Range("L141").Value = 73000
$ws.Range("N141").Value = -83360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 18171.477
$ws.Range("I3").Value = 22112.75
$ws.Range("K3").Value = 22112.75
$ws.Range("M3").Value = -21998.75

$ws.Range("H7").Value = 2776.3333
$ws.Range("I7").Value = 430
$ws.Range("K7").Value = 430
$ws.Range("M7").Value = -317

$ws.Range("H9").Value = 20000
$ws.Range("I9").Value = 20000
$ws.Range("K9").Value = 20000
$ws.Range("M9").Value = -19832

$ws.Range("H134").Value = 2373.818
$ws.Range("J134").Value = 3311.125
$ws.Range("L134").Value = 9933.375
$ws.Range("N134").Value = -15003.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1969.2858
$ws.Range("I31").Value = 1969.2858
$ws.Range("K31").Value = 1969.2858
$ws.Range("M31").Value = -1674.2858

$ws.Range("H34").Value = 1969.2858
$ws.Range("I34").Value = 1969.2858
$ws.Range("K34").Value = 1969.2858
$ws.Range("M34").Value = -1767.2858

$ws.Range("H54").Value = 25000
$ws.Range("J54").Value = 25000
$ws.Range("L54").Value = 25000
$ws.Range("N54").Value = -26316

$ws.Range("H109").Value = 25928.428
$ws.Range("J109").Value = 25928.428
$ws.Range("L109").Value = 25928.428
$ws.Range("N109").Value = -28008.428

$ws.Range("H134").Value = 4766.0835
$ws.Range("I134").Value = 4766.0835
$ws.Range("K134").Value = 14298.2505
$ws.Range("M134").Value = -11763.2505

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 99.5
$ws.Range("I44").Value = 99.5
$ws.Range("K44").Value = 298.5
$ws.Range("M44").Value = 99.5

$ws.Range("H118").Value = 5999.5713
$ws.Range("I118").Value = 5665.6665
$ws.Range("J118").Value = 6250
$ws.Range("K118").Value = 16996.9995
$ws.Range("L118").Value = 18750
$ws.Range("M118").Value = -15753.9995
$ws.Range("N118").Value = -21236

$ws.Range("H140").Value = 17859120
$ws.Range("I140").Value = 17859120
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 53577360
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -53572180
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 8135.3447
$ws.Range("I132").Value = 7237
$ws.Range("K132").Value = 21711
$ws.Range("M132").Value = -19181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 939.7143
$ws.Range("I55").Value = 916.6
$ws.Range("J55").Value = 997.5
$ws.Range("K55").Value = 916.6
$ws.Range("L55").Value = 997.5
$ws.Range("M55").Value = -743.6
$ws.Range("N55").Value = -1343.5

$ws.Range("H64").Value = 45000
$ws.Range("J64").Value = 45000
$ws.Range("L64").Value = 45000
$ws.Range("N64").Value = -45450

$ws.Range("H67").Value = 45000
$ws.Range("J67").Value = 45000
$ws.Range("L67").Value = 45000
$ws.Range("N67").Value = -46560

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 48383.848
$ws.Range("J115").Value = 48249.582
$ws.Range("L115").Value = 48249.582
$ws.Range("N115").Value = -51383.582

$ws.Range("H132").Value = 2723.3684
$ws.Range("I132").Value = 1771.2142
$ws.Range("J132").Value = 5389.4
$ws.Range("K132").Value = 5313.642599999999
$ws.Range("L132").Value = 16168.2
$ws.Range("M132").Value = -2783.642599999999
$ws.Range("N132").Value = -21228.2
